# Add json vlg, chelyaba, cheboksari
# Fix an existing Lada Largus model label, then append 17 new
# brand/model rows (UAZ, Lada, Mitsubishi, Haval, Suzuki, Great Wall,
# Dongfeng) to the catalog sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix existing row: "Фургон New" -> "Largus Фургон New" (Lada row 355) ---
$ws.Range("C355").Value = "Largus Фургон New"

# --- Append new catalog rows 586-602 ---
$newRows = @(
    @(685, "Mitsubishi", "ASX"),
    @(686, "Lada", "Largus CNG"),
    @(687, "Lada", "Largus Cross CNG"),
    @(688, "Lada", "Largus Фургон CNG"),
    @(689, "Haval", "M6 New"),
    @(690, "Suzuki", "Jimny"),
    @(691, "Great Wall", "GWM Wingle 7"),
    @(692, "Dongfeng", "DFM AX7"),
    @(693, "Lada", "Vesta Sedan New NG"),
    @(694, "UAZ", 3909),
    @(695, "UAZ", "3909 Бортовой грузовик"),
    @(696, "UAZ", "3909 Микроавтобус"),
    @(697, "UAZ", "3909 Скорая помощь"),
    @(698, "UAZ", "3909 Цельнометаллический фургон"),
    @(699, "UAZ", "Профи"),
    @(700, "Lada", "ВИС"),
    @(701, "Lada", "4x4 3 двери ")
)

$startRow = 586
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}

# --- Match the author's final selection/viewport (cosmetic) ---
$ws.Range("K604").Select()
